$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New column E: "a partir das 19:30" for rows 3-7 (reuses existing shared string)
$ws.Range("E3").Value = "a partir das 19:30"
$ws.Range("E4").Value = "a partir das 19:30"
$ws.Range("E5").Value = "a partir das 19:30"
$ws.Range("E6").Value = "a partir das 19:30"
$ws.Range("E7").Value = "a partir das 19:30"

# New column E: "dia todo" for rows 8-9 (new shared string)
$ws.Range("E8").Value = "dia todo"
$ws.Range("E9").Value = "dia todo"

# Give column E an explicit custom width (closest value reachable through
# the ColumnWidth property; the host quantizes to 1/6-character increments)
$ws.Columns(5).ColumnWidth = 17.8

# Update the selected / active cell to E6
$ws.Range("E6").Select() | Out-Null
